$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Datos actualizados" timestamp text in A1
$ws.Range("A1").Value = "Datos actualizados a 21 de Abril de 2020 a las 03:52"

# Row 27 - Corea del Sur
$ws.Range("B27").Value = 10683
$ws.Range("C27").Value = 9
$ws.Range("D27").Value = 8213
$ws.Range("E27").Value = 2233
$ws.Range("G27").Value = 1
$ws.Range("H27").Value = 237

# Row 49 - Panama
$ws.Range("B49").Value = 4658
$ws.Range("C49").Value = 191
$ws.Range("D49").Value = 204
$ws.Range("E49").Value = 4318
$ws.Range("G49").Value = 10
$ws.Range("H49").Value = 136

# Row 71 - Nueva Zelanda
$ws.Range("B71").Value = 1445
$ws.Range("C71").Value = 5
$ws.Range("D71").Value = 1006
$ws.Range("E71").Value = 426
$ws.Range("F71").Value = 3
$ws.Range("G71").Value = 1
$ws.Range("H71").Value = 13

# Row 117 - Guatemala
$ws.Range("B117").Value = 294
$ws.Range("C117").Value = 5
$ws.Range("D117").Value = 24
$ws.Range("E117").Value = 263
